{"js": "// Locate the \"A short report...\" bullet describing the project report\n// deliverable, rewrite its page-count guidance, and fold the separate\n// \"Graduate students...\" bullet into it (that whole paragraph goes away).\nconst body = context.document.body;\n\nconst oldText =\n  \"A short report describing your project and covering the details requested of your focus, no more than 1-2 pages.\";\nconst newText =\n  \"A short report describing your project and covering the details requested of your focus (3-6 pages, double spaced).\";\n\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target paragraph text not found\");\n}\n\nconst hitRange = results.items[0];\nconst targetPara = hitRange.paragraphs.getFirst();\n\n// The \"Graduate students...\" bullet immediately follows and is being\n// removed outright, so grab it before we touch the first paragraph.\nconst gradPara = targetPara.getNext();\ngradPara.load(\"text\");\nawait context.sync();\n\n// Replace the run text in the first bullet with the updated wording.\ntargetPara.insertText(newText, \"Replace\");\n\n// Drop the now-merged \"Graduate students...\" paragraph entirely.\nif (gradPara.text && gradPara.text.indexOf(\"Graduate students\") !== -1) {\n  gradPara.delete();\n} else {\n  // Fallback: locate it explicitly if paragraph ordering ever shifts.\n  const gradResults = body.search(\"Graduate students: this can be wrapped in with your 4-10 page paper.\", { matchCase: true });\n  gradResults.load(\"items\");\n  await context.sync();\n  if (gradResults.items.length > 0) {\n    gradResults.items[0].paragraphs.getFirst().delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the \"A short report...\" bullet's page-count guidance and fold the\n# separate \"Graduate students...\" bullet into it by removing that paragraph\n# outright.\n$d = $word.ActiveDocument\n\n$oldText = \"A short report describing your project and covering the details requested of your focus, no more than 1-2 pages.\"\n$newText = \"A short report describing your project and covering the details requested of your focus (3-6 pages, double spaced).\"\n$gradText = \"Graduate students: this can be wrapped in with your 4-10 page paper.\"\n\n# Locate the first bullet and rewrite its run text at the paragraph level\n# (rather than the narrower Find hit) so the existing run/paragraph\n# formatting is preserved instead of being replaced by a fresh run.\n$r = $d.Content\n$found = $r.Find.Execute($oldText)\nif ($found) {\n    $para = $r.Paragraphs(1)\n    $para.Range.Text = $newText\n}\n\n# Locate the now-orphaned \"Graduate students...\" paragraph and delete it\n# (including its paragraph mark) so the two bullets collapse into one.\n$r2 = $d.Content\n$found2 = $r2.Find.Execute($gradText)\nif ($found2) {\n    $gradPara = $r2.Paragraphs(1)\n    $gradPara.Range.Delete()\n}\n"}
